# Release update: bump the embedded build timestamp for
# "Coal Mine Boundaries and Methane Sources - version 1.0.0"
# from "February 03 2026 17.29.55 EST" to "February 03 2026 18.05.36 EST"
# everywhere it is referenced in the workbook.

$wb = $excel.ActiveWorkbook

$oldVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 18.05.36 EST)"

# --- "About" sheet ---
$wsAbout = $wb.Worksheets.Item("About")

# A2: "Version: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on ...)"
$a2 = $wsAbout.Range("A2")
$a2.Value = $a2.Value2.Replace($oldVersion, $newVersion)

# A6: "Recommended Citation: ... version '<version>'. (See the CC license ...)"
$a6 = $wsAbout.Range("A6")
$a6.Value = $a6.Value2.Replace($oldVersion, $newVersion)

# --- "Boundaries and methane sources" sheet ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")
$usedRange = $wsData.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $wsData.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null -and $val -eq $oldVersion) {
            $cell.Value = $newVersion
        }
    }
}
